# Challenge Problem 6.pptx edit
#
# 1) Slide 18 ("CP#6 Next Evaluation Period"), Content Placeholder:
#    - merge the two runs "October 15-30: Beta " + "Period" into one run
#    - merge the three runs "45 Days before PI meeting: Final " +
#      "Deadline for CP5 and CP6 " + "solutions" into one run
#    - merge the two runs "co-located " + "but " into one run
# 2) Slide 19 ("CP#6 Materials Available Now", final slide), Content
#    Placeholder: append a blank line, an "Email address for questions,
#    issues, etc.:" line, and a line with a tab + a mailto: hyperlink to
#    ppaml-support@community.galois.com, then a trailing blank line.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Part 1: slide 18 - collapse split runs that share identical formatting
# ---------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$shp18 = $s18.Shapes.Item(2)   # "Content Placeholder 3"
$tr18 = $shp18.TextFrame.TextRange

# Paragraph 2: "October 15-30: Beta " + "Period" -> "October 15-30: Beta Period"
$para2 = $tr18.Paragraphs(2, 1)
$para2.Text = "TEMP_MERGE_PLACEHOLDER_1"
$para2 = $tr18.Paragraphs(2, 1)
$para2.Text = "October 15-30: Beta Period"

# Paragraph 3: "45 Days before PI meeting: Final " + "Deadline for CP5 and
# CP6 " + "solutions" -> single run
$para3 = $tr18.Paragraphs(3, 1)
$para3.Text = "TEMP_MERGE_PLACEHOLDER_2"
$para3 = $tr18.Paragraphs(3, 1)
$para3.Text = "45 Days before PI meeting: Final Deadline for CP5 and CP6 solutions"

# Paragraph 4: only merge the inner "co-located " + "but " runs, leaving
# the surrounding "January ??: PI Meeting (" and "not conflicting with
# POPL)" runs (which carry different rPr) untouched.
$fullText18 = $tr18.Text
$idx0 = $fullText18.IndexOf("co-located ")
$startPos = $idx0 + 1
$oldLen = ("co-located " + "but ").Length
$sub = $tr18.Characters($startPos, $oldLen)
$sub.Text = "TEMP_MERGE_PLACEHOLDER_3"
$sub = $tr18.Characters($startPos, 24)
$sub.Text = "co-located but "

# ---------------------------------------------------------------------
# Part 2: slide 19 (last slide) - add the email contact block
# ---------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$shp19 = $s19.Shapes.Item(2)   # "Content Placeholder 3"
$tr19 = $shp19.TextFrame.TextRange

$tab = [char]9
$newBlock = "`r`rEmail address for questions, issues, etc.:`r" + $tab + "ppaml-support@community.galois.com`r"
$tr19.InsertAfter($newBlock)

$count19 = $tr19.Paragraphs().Count

# Paragraph holding the tab + email address is the second-to-last new
# paragraph; the trailing blank paragraph is last.
$emailPara = $tr19.Paragraphs($count19 - 1, 1)
$introPara = $tr19.Paragraphs($count19 - 2, 1)
$blankPara = $tr19.Paragraphs($count19 - 3, 1)
$trailingPara = $tr19.Paragraphs($count19, 1)

# Remove the default outline bullet / indent from the three new text
# paragraphs (blank separator keeps the outline bullet format from the
# body placeholder, matching the un-bulleted "buNone" paragraphs added
# in the source edit).
$introPara.ParagraphFormat.Bullet.Type = 0
$emailPara.ParagraphFormat.Bullet.Type = 0
$trailingPara.ParagraphFormat.Bullet.Type = 0

# Split the tab + email-address paragraph into its own two runs so the
# hyperlink/underline formatting only applies to the address itself.
$emailFullText = $emailPara.Text
$emailStart = $emailPara.Start + 1  # 1-based index of the character right after the tab
$addr = $tr19.Characters($emailStart, ("ppaml-support@community.galois.com").Length)
$addrAction = $addr.ActionSettings(1)
$addrAction.Hyperlink.Address = "mailto:ppaml-support@community.galois.com"
$addr.Font.Underline = -1

Write-Host "Slide 18 content:" $tr18.Text
Write-Host "Slide 19 content:" $tr19.Text
